# Add new columns I (I0) and J (IF) with header style matching existing headers (style index 1),
# and fill in numeric data for rows 2-12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - copy style from an existing header cell (H1) so it matches (bold, centered, bordered)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-12 for columns I and J
$data = @(
    @(2, 3, 4),
    @(3, 9, 9),
    @(4, 1, 2),
    @(5, 5, 8),
    @(6, 7, 8),
    @(7, 5, 6),
    @(8, 6, 8),
    @(9, 9, 9),
    @(10, 1, 4),
    @(11, 1, 3),
    @(12, 1, 2)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($row, 9).Value = $iVal
    $ws.Cells.Item($row, 10).Value = $jVal
}
